$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ===========================================================================
# Block 1: "General X8" merged group (A14:A24 -> A14:A26).
# Two brand-new rows (5.1.9.1 / 5.1.9.0) are inserted above the existing
# 5.1.8.1 entry, which simply shifts down along with the rest of the block.
# ===========================================================================

$generalLabel = $ws.Range("A14").Value2

$ws.Rows("14:15").Insert()

$ws.Range("A17:D17").Copy()
$ws.Range("A14:D15").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A16").Value = ""
$ws.Range("A14").Value = $generalLabel
$ws.Range("A14:A26").Merge()
$ws.Rows("14:15").RowHeight = 15.75

$ws.Range("B14").Value = "5.1.9.1"
$ws.Range("C14").Value = "X8-A103/AIBT"
$ws.Range("D14").Value = "Features required for version 5.1.9.0"

$ws.Range("B15").Value = "5.1.9.0"
$ws.Range("C15").Value = "X8-A107/F21"
$ws.Range("D15").Value = "XYR coordinate conversion table add export xy function #395"

# Comment on the (now shifted) B16 cell follows the original B14 cell.
$cmt1 = $ws.Range("B14").Comment
$cmtText1 = $cmt1.Text()
$cmt1.Delete()
$null = $ws.Range("B16").AddComment($cmtText1)

# ===========================================================================
# Block 2: "Previous X8 / X6" merged group (A29:A33 -> A31:A37, after the
# +2 shift from block 1 this group now starts at row 31).
# Two brand-new rows (5.4.5.3 / 5.4.5.2) are inserted above the existing
# 5.4.5.1 entry.
# ===========================================================================

$prevLabel = $ws.Range("A31").Value2

$ws.Rows("31:32").Insert()

$ws.Range("A34:D34").Copy()
$ws.Range("A31:D32").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A33").Value = ""
$ws.Range("A31").Value = $prevLabel
$ws.Range("A31:A37").Merge()
$ws.Rows("31:32").RowHeight = 15.75

$ws.Range("B31").Value = "5.4.5.3"
$ws.Range("C31").Value = "X6-B102/F12A P4 EBO"
$ws.Range("D31").Value = "For new GV2 (hobochen)"

$ws.Range("B32").Value = "5.4.5.2"
$ws.Range("C32").Value = "X6-B102/F12A P4 EBO"
$ws.Range("D32").Value = "For new GV2 (temiceng)"

# Comment on the (now shifted) B33 cell follows the original B29 cell.
$cmt2 = $ws.Range("B31").Comment
$cmtText2 = $cmt2.Text()
$cmt2.Delete()
$null = $ws.Range("B33").AddComment($cmtText2)

# ===========================================================================
# Misc: the red highlight font used for the most-recently-shipped version in
# each block is recoloured from FF0000 to C00000.
# ===========================================================================
$ws.Range("B16:D16").Font.Color = 192   # RGB(192,0,0) == C00000, COM BGR long
$ws.Range("B33:D33").Font.Color = 192

Write-Host "All steps complete"
